$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.343.33'
$ws.Range('E2').Value = '  -0.42%  '

$ws.Range('D3').Value = '1.714.39'
$ws.Range('E3').Value = '  -0.36%  '

$ws.Range('E4').Value = '  +0.28%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '224.74'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.03%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5278'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.04%  '

$ws.Range('E7').Value = '  +0.26%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.06666'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.22%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2644'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.38%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.76'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.24%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07755'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.15%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.468'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.89%  '

$ws.Range('D13').Value = '1.950.61'
$ws.Range('E13').Value = '  -0.43%  '

$ws.Range('D14').Value = '1.719.50'
$ws.Range('E14').Value = '  -0.20%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5793'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.31%  '

$ws.Range('D16').Value = '0.0₅8176'
$ws.Range('E16').Value = '  -0.89%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '67.71'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.06%  '

$ws.Range('D18').Value = '27.356.45'
$ws.Range('E18').Value = '  -0.37%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '219.54'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.81%  '

$ws.Range('E20').Value = '  +0.37%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.649'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.48%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.42'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.13%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.033'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.89%  '

$ws.Range('E24').Value = '  +0.27%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.16'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.45%  '

$ws.Range('E26').Value = '  -1.59%  '

$ws.Range('E27').Value = '  -1.76%  '

$ws.Range('E28').Value = '  -1.07%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '16.17'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.57%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05354'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.41%  '

$ws.Range('E31').Value = '  -0.19%  '

$ws.Range('E32').Value = '  -1.25%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.388'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.17%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.636'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.70%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.847'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.15%  '

$ws.Range('E36').Value = '  -0.22%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.401'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.08%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5876'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.30%  '

$ws.Range('D39').Value = '1.159.15'
$ws.Range('E39').Value = '  +11.03%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01649'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.62%  '

$ws.Range('E41').Value = '  -1.26%  '

$ws.Range('E42').Value = '  +0.29%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8397'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.63%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '101.12'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.02%  '

$ws.Range('D45').Value = '1.857.07'
$ws.Range('E45').Value = '  -0.49%  '

$ws.Range('E46').Value = '  +1.37%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '57.40'
$ws.Range('D47').Style = "Normal"

$ws.Range('E48').Value = '  +0.83%  '

$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.005'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.01%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.130'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.78%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05242'
$ws.Range('D51').Style = "Normal"
